# WMS 11866 - Data Platform Workshop 2025-04-28 update
#
# 1) Bump the cached "datetimeFigureOut" footer-date field from 4/23/25 to
#    4/28/25 everywhere it is stamped: the slide master's Date Placeholder
#    and the Date Placeholder on every slide layout.
# 2) Give the "Rectangle 3" title-bar shape (inside "Group 4" on slide 1)
#    an explicit "no line" outline.

$p = $ppt.ActivePresentation

$oldDate = "4/23/25"
$newDate = "4/28/25"

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's own Date Placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's Date Placeholder.
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 1: "Rectangle 3" (child of group "Group 4") gets an explicit
# no-line outline alongside its existing solid 153040 fill.
$slide = $p.Slides.Item(1)
for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
    $topShape = $slide.Shapes.Item($k)
    if ($topShape.Name -eq "Group 4") {
        $items = $topShape.GroupItems
        for ($m = 1; $m -le $items.Count; $m++) {
            $child = $items.Item($m)
            if ($child.Name -eq "Rectangle 3") {
                $child.Line.Visible = $false
            }
        }
    }
}
